$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.167.81"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.827.34"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "1.821.93"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "29.145.31"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "2.067.91"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  +6.94%  "
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  +3.40%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("E38").Value = "  +2.88%  "
$ws.Range("D39").Value = "1.214.15"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "1.975.56"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("E51").Value = "  -1.21%  "

# Cells whose new values look numeric but must remain exact text;
# force text format, assign, then restore default style so no stray
# cell-level style attribute is left behind.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6223"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07347"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2903"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07685"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.966"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6667"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008975"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.859"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.360"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1425"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.493"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05573"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.093"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.100"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.205"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.846"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7358"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.844"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01769"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.309"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9137"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5089"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000118"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.144"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4032"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05761"
$ws.Range("D51").Style = "Normal"
